$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D) and Volume(1h) (E) values for rows with Price changes ---
# D-column cells are text (they use "." as a thousands separator, e.g. "43.033.91")
# so NumberFormat is forced to Text before the write, then restored to General,
# to avoid Excel auto-converting the numeric-looking strings into real numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.033.91"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.305.89"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.66"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -0.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.13"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -1.66%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.511"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -1.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.37"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -0.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0789"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  -0.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "17.91"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  +0.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.81"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -1.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.664.25"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -0.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.314.15"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -2.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.784"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  -2.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.005.52"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.72"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +0.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "242.36"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +1.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.14"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -1.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.17"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -1.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.63"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -1.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.31"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -3.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.75"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +0.24%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.75"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +0.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0689"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -1.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.77"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +0.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.998.34"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -0.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.24"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +1.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.43"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -1.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.80"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -3.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.66"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -2.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.530.34"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -0.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.90"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -1.61%  "

# --- Update Volume(1h) (E) only values ---
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("E7").Value = "  +2.34%  "
$ws.Range("E13").Value = "  +0.63%  "
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("E21").Value = "  -1.16%  "
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("E26").Value = "  -1.16%  "
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("E34").Value = "  -3.56%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("E39").Value = "  -1.11%  "
$ws.Range("E40").Value = "  -1.79%  "
$ws.Range("E42").Value = "  +0.86%  "
$ws.Range("E44").Value = "  -1.63%  "
$ws.Range("E45").Value = "  -2.86%  "

# --- Rows 30 and 31 swap: Toncoin moves to row 30, Cosmos moves to row 31 ---
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.04"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -0.41%  "
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.10"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -0.86%  "
